# Update the "Expected Wins" sheet's Difference column (D2:D13) so it
# reflects Expected Wins (C) minus actual Wins (parsed from the Record
# text in column E), instead of the previous (incorrect) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expected Wins")

$newValues = @{
    2  = 0.333333333333333
    3  = 0.25
    4  = 0
    5  = -1.333333333333333
    6  = -0.5
    7  = 0.3333333333333335
    8  = -0.08333333333333348
    9  = 1.916666666666667
    10 = -0.3333333333333335
    11 = -0.9166666666666665
    12 = 0.4166666666666667
    13 = -0.08333333333333337
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $newValues[$row]
}
